# illuminance vs angle.xlsx -- "report-poster: Add poster and presentation files"
#
# Summary of the edit being replayed:
#  1. A new worksheet ("Sheet11") is appended after "Sheet10". It collects the
#     0-degree ('100% brightness'!F8:F16) illuminance readings into a small
#     -90..90 degree table (mirrors the layout already used on "Blue").
#  2. The "Blue" worksheet gets a mirrored copy of its C:G illuminance table
#     in columns K:Q (symmetric around 0 degrees: K/Q=G, L/P=F, M/O=E, N=D)
#     for rows 8-17, and that range becomes the active selection.
#  3. The scatter/trendline chart embedded on "distance and illuminance" is
#     removed.
#  4. "Blue" becomes the active sheet/tab (previously "distance and
#     illuminance" was active).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add "Sheet11" after the last existing sheet ("Sheet10").
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet11 = $wb.Worksheets.Add($null, $lastSheet)
$sheet11.Name = "Sheet11"

$sheet11.Range("B6").Value = 20

$sheet11.Range("C36").Value = -90
$sheet11.Range("D36").Value = -60
$sheet11.Range("E36").Value = -30
$sheet11.Range("F36").Value = 0
$sheet11.Range("G36").Value = 30
$sheet11.Range("H36").Value = 60
$sheet11.Range("I36").Value = 90

$sheet11.Range("B37").Value = 20
$sheet11.Range("B38").Value = 30
$sheet11.Range("B39").Value = 40
$sheet11.Range("B40").Value = 50
$sheet11.Range("B41").Value = 60
$sheet11.Range("B42").Value = 70
$sheet11.Range("B43").Value = 80
$sheet11.Range("B44").Value = 90
$sheet11.Range("B45").Value = 100

$sheet11.Range("I37").Formula = "='100% brightness'!F8"
$sheet11.Range("I38").Formula = "='100% brightness'!F9"
$sheet11.Range("I39").Formula = "='100% brightness'!F10"
$sheet11.Range("I40").Formula = "='100% brightness'!F11"
$sheet11.Range("I41").Formula = "='100% brightness'!F12"
$sheet11.Range("I42").Formula = "='100% brightness'!F13"
$sheet11.Range("I43").Formula = "='100% brightness'!F14"
$sheet11.Range("I44").Formula = "='100% brightness'!F15"
$sheet11.Range("I45").Formula = "='100% brightness'!F16"

# the new style used for these cells bumps the shared "Arial/black" font
# from 10pt to 11pt
$sheet11.Range("I37:I45").Font.Size = 11

# ---------------------------------------------------------------------------
# 2. Mirror the Blue sheet's C:G illuminance table into K:Q (rows 8-17) so
#    the data spans the full -90..90 degree range.
# ---------------------------------------------------------------------------
$blue = $wb.Worksheets.Item("Blue")

for ($row = 8; $row -le 17; $row++) {
    $blue.Range("K" + $row).Formula = "=G" + $row
    $blue.Range("L" + $row).Formula = "=F" + $row
    $blue.Range("M" + $row).Formula = "=E" + $row
    $blue.Range("N" + $row).Formula = "=D" + $row
    $blue.Range("O" + $row).Formula = "=E" + $row
    $blue.Range("P" + $row).Formula = "=F" + $row
    $blue.Range("Q" + $row).Formula = "=G" + $row
}

# ---------------------------------------------------------------------------
# 3. Remove the scatter/trendline chart from "distance and illuminance".
# ---------------------------------------------------------------------------
$dist = $wb.Worksheets.Item("distance and illuminance")
$charts = $dist.ChartObjects()
for ($i = $charts.Count; $i -ge 1; $i--) {
    $charts.Item($i).Delete()
}

# ---------------------------------------------------------------------------
# 4. "Blue" becomes the active sheet with K8:Q17 selected (matches the new
#    activeTab index and the sheetView selection recorded for Blue).
# ---------------------------------------------------------------------------
$blue.Activate()
$blue.Range("K8:Q17").Select()
